$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '''65.934.01'
$ws.Range("E2").Value = '  +0.48%  '
$ws.Range("D3").Value = '''2.666.80'
$ws.Range("E3").Value = '  -0.02%  '
$ws.Range("D4").Value = '''1.00'
$ws.Range("E4").Value = '  +0.12%  '
$ws.Range("D5").Value = '''600.35'
$ws.Range("E5").Value = '  +0.27%  '
$ws.Range("D6").Value = '''161.07'
$ws.Range("E6").Value = '  +2.94%  '
$ws.Range("D7").Value = '''0.644'
$ws.Range("E7").Value = '  +4.48%  '
$ws.Range("D8").Value = '''1.00'
$ws.Range("E8").Value = '  +0.04%  '
$ws.Range("D9").Value = '''0.127'
$ws.Range("E9").Value = '  -1.10%  '
$ws.Range("E10").Value = '  +0.93%  '
$ws.Range("E11").Value = '  +0.92%  '
$ws.Range("E12").Value = '  +1.67%  '
$ws.Range("D13").Value = '''29.38'
$ws.Range("E13").Value = '  +0.45%  '
$ws.Range("D14").Value = '''0.0000196'
$ws.Range("E14").Value = '  +0.22%  '
$ws.Range("D15").Value = '''3.147.02'
$ws.Range("E15").Value = '  +0.02%  '
$ws.Range("D16").Value = '''65.810.18'
$ws.Range("E16").Value = '  +0.56%  '
$ws.Range("D17").Value = '''2.641.03'
$ws.Range("E17").Value = '  -0.87%  '
$ws.Range("E18").Value = '  -0.92%  '
$ws.Range("D19").Value = '''4.83'
$ws.Range("E19").Value = '  +1.32%  '
$ws.Range("D20").Value = '''356.91'
$ws.Range("E20").Value = '  +1.67%  '
$ws.Range("D21").Value = '''7.51'
$ws.Range("E21").Value = '  -0.38%  '
$ws.Range("E22").Value = '  -0.04%  '
$ws.Range("D23").Value = '''70.19'
$ws.Range("E23").Value = '  +0.81%  '
$ws.Range("D24").Value = '''1.80'
$ws.Range("E24").Value = '  +10.26%  '
$ws.Range("D25").Value = '''0.0000114'
$ws.Range("E25").Value = '  +3.17%  '
$ws.Range("D26").Value = '''9.80'
$ws.Range("E26").Value = '  +2.27%  '
$ws.Range("E27").Value = '  +2.81%  '
$ws.Range("D28").Value = '''579.97'
$ws.Range("E28").Value = '  +11.01%  '
$ws.Range("E29").Value = '  +2.47%  '
$ws.Range("E30").Value = '  -1.52%  '
$ws.Range("B31").Value = 'PancakeSwap'
$ws.Range("C31").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D31").Value = '''2.16'
$ws.Range("E31").Value = '  +1.58%  '
$ws.Range("B32").Value = 'Binance-PegBSC-USD'
$ws.Range("C32").Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Range("D32").Value = '''0.999'
$ws.Range("E32").Value = '  -0.14%  '
$ws.Range("D33").Value = '''1.84'
$ws.Range("E33").Value = '  +4.62%  '
$ws.Range("D34").Value = '''6.76'
$ws.Range("E34").Value = '  +5.22%  '
$ws.Range("D35").Value = '''5.52'
$ws.Range("E35").Value = '  +1.37%  '
$ws.Range("E36").Value = '  +0.71%  '
$ws.Range("D37").Value = '''20.67'
$ws.Range("E37").Value = '  +0.66%  '
$ws.Range("E38").Value = '  +2.65%  '
$ws.Range("D39").Value = '''1.00'
$ws.Range("E39").Value = '  +0.15%  '
$ws.Range("D40").Value = '''154.72'
$ws.Range("E40").Value = '  -1.84%  '
$ws.Range("D41").Value = '''2.51'
$ws.Range("E41").Value = '  +9.39%  '
$ws.Range("D42").Value = '''162.56'
$ws.Range("E42").Value = '  -0.28%  '
$ws.Range("E43").Value = '  +0.54%  '
$ws.Range("E44").Value = '  +2.59%  '
$ws.Range("D45").Value = '''23.65'
$ws.Range("E45").Value = '  +4.30%  '
$ws.Range("E46").Value = '  +1.53%  '
$ws.Range("E47").Value = '  +1.33%  '
$ws.Range("E48").Value = '  +2.19%  '
$ws.Range("D49").Value = '''19.89'
$ws.Range("E49").Value = '  -0.77%  '
$ws.Range("E50").Value = '  -5.48%  '
$ws.Range("D51").Value = '''0.823'
$ws.Range("E51").Value = '  +1.97%  '
